$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 16096.286
$ws.Range("I70").Value = 880
$ws.Range("K70").Value = 2640
$ws.Range("M70").Value = -2370
$ws.Range("H73").Value = 16096.286
$ws.Range("I73").Value = 880
$ws.Range("K73").Value = 2640
$ws.Range("M73").Value = -1704
$ws.Range("H100").Value = 6431.5264
$ws.Range("J100").Value = 7664.7
$ws.Range("L100").Value = 7664.7
$ws.Range("N100").Value = -8746.700000000001
$ws.Range("H107").Value = 1527.8889
$ws.Range("I107").Value = 2008.5834
$ws.Range("K107").Value = 2008.5834
$ws.Range("M107").Value = -88.58339999999998
$ws.Range("H131").Value = 13250.875
$ws.Range("I131").Value = 9668.223
$ws.Range("J131").Value = 17857.143
$ws.Range("K131").Value = 29004.669
$ws.Range("L131").Value = 53571.429
$ws.Range("M131").Value = -23964.669
$ws.Range("N131").Value = -63651.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4995.6577
$ws.Range("I32").Value = 5121.1113
$ws.Range("J32").Value = 2737.5
$ws.Range("K32").Value = 5121.1113
$ws.Range("L32").Value = 2737.5
$ws.Range("M32").Value = -4834.1113
$ws.Range("N32").Value = -3311.5
$ws.Range("H74").Value = 1832.75
$ws.Range("I74").Value = 1957.5238
$ws.Range("K74").Value = 1957.5238
$ws.Range("M74").Value = -1083.5238
$ws.Range("H77").Value = 1832.75
$ws.Range("I77").Value = 1957.5238
$ws.Range("K77").Value = 9787.618999999999
$ws.Range("M77").Value = -5419.618999999999
$ws.Range("H110").Value = 3909.818
$ws.Range("I110").Value = 3714
$ws.Range("J110").Value = 4575.6
$ws.Range("K110").Value = 3714
$ws.Range("L110").Value = 4575.6
$ws.Range("M110").Value = -1669
$ws.Range("N110").Value = -8665.6
$ws.Range("H132").Value = 1875.069
$ws.Range("I132").Value = 2017.9524
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 6053.857199999999
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -3523.857199999999
$ws.Range("N132").Value = -9560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 304.25
$ws.Range("J80").Value = 291.42856
$ws.Range("L80").Value = 291.42856
$ws.Range("N80").Value = -2287.42856
$ws.Range("H83").Value = 304.25
$ws.Range("J83").Value = 291.42856
$ws.Range("L83").Value = 1457.1428
$ws.Range("N83").Value = -11441.1428
$ws.Range("H134").Value = 3270.7144
$ws.Range("I134").Value = 2644.111
$ws.Range("J134").Value = 4398.6
$ws.Range("K134").Value = 7932.333
$ws.Range("L134").Value = 13195.8
$ws.Range("M134").Value = -5397.333
$ws.Range("N134").Value = -18265.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3337.8823
$ws.Range("I16").Value = 3188.1667
$ws.Range("K16").Value = 3188.1667
$ws.Range("M16").Value = -2901.1667
$ws.Range("H18").Value = 57500
$ws.Range("J18").Value = 57500
$ws.Range("L18").Value = 57500
$ws.Range("N18").Value = -57960
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51498
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -157488
$ws.Range("H105").Value = 626.4167
$ws.Range("I105").Value = 592.4545000000001
$ws.Range("K105").Value = 592.4545000000001
$ws.Range("M105").Value = 1154.5455
$ws.Range("H107").Value = 8935865
$ws.Range("J107").Value = 17820
$ws.Range("L107").Value = 17820
$ws.Range("N107").Value = -21660
$ws.Range("H113").Value = 3337.8823
$ws.Range("I113").Value = 3188.1667
$ws.Range("K113").Value = 3188.1667
$ws.Range("M113").Value = -1018.1667
$ws.Range("H132").Value = 2663.2812
$ws.Range("I132").Value = 2037.4445
$ws.Range("K132").Value = 6112.333500000001
$ws.Range("M132").Value = -3582.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 61.5
$ws.Range("I2").Value = 58.25
$ws.Range("K2").Value = 349.5
$ws.Range("M2").Value = -236.5
$ws.Range("H5").Value = 799.15625
$ws.Range("J5").Value = 1134.25
$ws.Range("L5").Value = 3402.75
$ws.Range("N5").Value = -3626.75
$ws.Range("H131").Value = 3482.25
$ws.Range("I131").Value = 1842.2222
$ws.Range("K131").Value = 5526.6666
$ws.Range("M131").Value = -486.6665999999996
$ws.Range("H135").Value = 799.15625
$ws.Range("J135").Value = 1134.25
$ws.Range("L135").Value = 10208.25
$ws.Range("N135").Value = -15278.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 105228.91
$ws.Range("J70").Value = 5310.8
$ws.Range("L70").Value = 5310.8
$ws.Range("N70").Value = -5850.8
$ws.Range("H73").Value = 105228.91
$ws.Range("J73").Value = 5310.8
$ws.Range("L73").Value = 5310.8
$ws.Range("N73").Value = -7182.8
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("H107").Value = 652.2273
$ws.Range("I107").Value = 551.8125
$ws.Range("K107").Value = 551.8125
$ws.Range("M107").Value = 1368.1875
$ws.Range("H113").Value = 5915.8
$ws.Range("I113").Value = 2745.7778
$ws.Range("K113").Value = 2745.7778
$ws.Range("M113").Value = -575.7777999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6889.6665
$ws.Range("I7").Value = 2669
$ws.Range("J7").Value = 9000
$ws.Range("K7").Value = 2669
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = -2557
$ws.Range("N7").Value = -9224
$ws.Range("H126").Value = 6889.6665
$ws.Range("I126").Value = 2669
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 8007
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -5537
$ws.Range("N126").Value = -31940
$ws.Range("H132").Value = 1647.6428
$ws.Range("I132").Value = 1511.1316
$ws.Range("J132").Value = 2944.5
$ws.Range("K132").Value = 4533.3948
$ws.Range("L132").Value = 8833.5
$ws.Range("M132").Value = -2003.3948
$ws.Range("N132").Value = -13893.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2542.6365
$ws.Range("I81").Value = 1784.7142
$ws.Range("J81").Value = 3869
$ws.Range("K81").Value = 3569.4284
$ws.Range("L81").Value = 7738
$ws.Range("M81").Value = -2508.4284
$ws.Range("N81").Value = -9860
$ws.Range("H84").Value = 2542.6365
$ws.Range("I84").Value = 1784.7142
$ws.Range("J84").Value = 3869
$ws.Range("K84").Value = 17847.142
$ws.Range("L84").Value = 38690
$ws.Range("M84").Value = -12543.142
